# Save changes before implementing dynamic form customization
#
# The sheet previously tracked warranty/AMC fields (with 2 sample asset
# rows). It is being repurposed towards a simpler "amc_contract / end_user"
# layout with a single in-progress data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -------------------------------------------------------
# Columns A:E (serial_number, asset_type, po_number, sap_asset_id,
# installation_date) are unchanged. F/G are repurposed; H:K are dropped.
$ws.Range("F1").Value = "amc_contract"
$ws.Range("G1").Value = "end_user"
$ws.Range("H1:K1").Clear()

# --- Data row 2 ---------------------------------------------------------
# New single data row: serial number + asset type, plus an amc contract
# number and the end user's name. The old PO/SAP/date/provider columns
# for this row are cleared.
$ws.Range("A2").Value = 12345
$ws.Range("B2").Value = "laptop"
$ws.Range("C2:E2").Clear()
$ws.Range("F2").Value = 611888
$ws.Range("G2").Value = "Tituraj"
$ws.Range("H2:K2").Clear()

# --- Old second data row (row 3) removed entirely ------------------------
$ws.Range("A3:K3").Clear()

# --- Column widths --------------------------------------------------------
# Column A now holds the serial number and gets an explicit best-fit-style
# width; column B keeps its existing width untouched.
$ws.Columns.Item(1).ColumnWidth = 12.14

# --- Selection ------------------------------------------------------------
$ws.Range("F3").Select()
